# Atualizacao rapida de agenda as  8:48:10,92
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Pedro / Quinta Avenida - porta de vidro (H2/I2 untouched)
$ws.Range("A2").Value = "Pedro"
$ws.Range("B2").Value = "0007"
$ws.Range("C2").Value = "Quinta Avenida"
$ws.Range("D2").Value = "Marcos pediu pra ir lá, acho que tem haver com a porta de vidro."
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Pendente"

# Row 3: Pedro / Wetravel - central off
$ws.Range("A3").Value = "Pedro"
$ws.Range("B3").Value = "2090"
$ws.Range("C3").Value = "Wetravel"
$ws.Range("D3").Value = "Central off, cliente disse que mudou de internet."
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "Pendente"
$ws.Rows.Item(3).RowHeight = 15

# Row 4: Roberto / Gerdau Galpão - revisão (continuação)
$ws.Range("A4").Value = "Roberto"
$ws.Range("B4").Value = "0638"
$ws.Range("C4").Value = "Gerdau Galpão"
$ws.Range("D4").Value = "Revisão no sistema do cliente, continuação de ontem."
$ws.Range("E4").Value = ""
$ws.Range("G4").Value = "Em andamento"

# Row 5: Roberto / Gerdau Escritório - revisão (continuação)
$ws.Range("A5").Value = "Roberto"
$ws.Range("B5").Value = "0576"
$ws.Range("C5").Value = "Gerdau Escritório"
$ws.Range("D5").Value = "Revisão no sistema do cliente, continuação de ontem."
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "Em andamento"

# Row 6: Giovani / MegaScan - zona aberta
$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "0643"
$ws.Range("C6").Value = "MegaScan"
$ws.Range("D6").Value = "Zona aberta, cliente pedindo reparo."
$ws.Range("E6").Value = ""
$ws.Range("G6").Value = "Em andamento"

# Rows 7-11: fully cleared (old agenda entries removed)
$ws.Range("A7:I11").ClearContents
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15

# Selection moves to H6
$ws.Range("H6").Select
